$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 8 data: CR 8.5 Person Search Page Changes ---
# Column mapping: A=Date, B=Bugzilla Ticket No, C=Artifacts Name, D=Comments,
#                 E=Developer, F=Require Database, G=SQL

# Set values in this order so new shared strings append as 29, 30, 31
$ws.Range("D8").Value = "CR 8.5 Person Search Page Changes"
$ws.Range("F8").Value = "No"
$longText = @"
New Files Added:
\webclient\components\MOLSA\ReferenceApp\core\Person\Search\MOLSAPerson_searchCriteriaView.properties
\webclient\components\MOLSA\ReferenceApp\core\Person\Search\MOLSAPerson_searchCriteriaView.vim
\webclient\components\MOLSA_ar\ReferenceApp\core\Person\Search\MOLSAPerson_searchCriteriaView_ar.properties
\EJBServer\components\MOLSA\message\MOLSABpoPersonSearch.xml
\EJBServer\components\MOLSA\source\curam\molsa\core\impl\MOLSADatabasePersonSearch.java 
Existing Files Mofified:
\webclient\components\MOLSA\ReferenceApp\core\Person\Search\Person_search1.uim  
\EJBServer\components\MOLSA\model\Packages\Reference Model\Core.efx  
\EJBServer\components\MOLSA\model\Packages\Reference Model\Facade\Person.efx  
\EJBServer\components\MOLSA\model\Packages\Reference Model\Facade.efx  
\EJBServer\components\MOLSA\model\Packages\Reference Model\Person Search\Person Search.efx  
\EJBServer\components\MOLSA\source\curam\molsa\core\facade\impl\MOLSAPersonDA.java  
"@
$ws.Range("C8").Value = $longText

$ws.Range("E8").Value = "Bhaskar"

# Date: March 17, 2017 (serial 42811) - set as raw serial so it inherits
# the existing column A date style instead of creating a new number format
$ws.Range("A8").Value = 42811

# D8 needs word-wrap like the other "Comments" cells
$ws.Range("D8").WrapText = $true

# --- Row heights (content reflow after the new row) ---
$ws.Rows.Item(1).RowHeight = 43.2
$ws.Rows.Item(3).RowHeight = 124.2
$ws.Rows.Item(4).RowHeight = 408.6
$ws.Rows.Item(6).RowHeight = 129.6
$ws.Rows.Item(7).RowHeight = 360
$ws.Rows.Item(8).RowHeight = 216

# --- Column widths (minor re-fit) ---
$ws.Columns.Item(1).ColumnWidth = 17.053385416666668
$ws.Columns.Item(4).ColumnWidth = 34.608072916666664
$ws.Columns.Item(5).ColumnWidth = 12.053385416666666
$ws.Columns.Item(6).ColumnWidth = 11.721354166666666
$ws.Columns.Item(7).ColumnWidth = 187.05338541666666
$ws.Columns.Item(8).ColumnWidth = 11.276041666666666
$ws.Columns.Item(9).ColumnWidth = 52.276041666666664

# --- Selection moves to the newly-entered row ---
$ws.Range("A8").Select()
